$wb = $excel.ActiveWorkbook

# --- Data fix: stationsExperiment!H7 410 -> 350 ---
$wsStationsExperiment = $wb.Worksheets.Item("stationsExperiment")
$wsStationsExperiment.Range("H7").Value = 350

# --- Per-sheet selection updates (each worksheet remembers its own last selection) ---

# Notifications: B10 -> D39
$wsNotifications = $wb.Worksheets.Item("Notifications")
$wsNotifications.Activate()
$wsNotifications.Range("D39").Select()

# stationsLearn: G2:G4 -> C39 (single cell)
$wsStationsLearn = $wb.Worksheets.Item("stationsLearn")
$wsStationsLearn.Activate()
$wsStationsLearn.Range("C39").Select()

# interruptionsExperiment: N3:N19 -> I9 (single cell)
$wsInterruptionsExperiment = $wb.Worksheets.Item("interruptionsExperiment")
$wsInterruptionsExperiment.Activate()
$wsInterruptionsExperiment.Range("I9").Select()

# stationsExperiment: G4:G7 -> H8 (single cell); becomes the active tab/window scroll
$wsStationsExperiment.Activate()
$wsStationsExperiment.Range("H8").Select()

# --- Window view: scroll tabs so Setup (index 1) is the first visible tab, window at origin ---
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.ScrollWorkbookTabs(1)
